# Add Files Upload and Update Data
# Adds a "Path" column (C) with image-path lookups for each employee row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Path"

# New path values per row - leading "'" reproduces Excel's automatic
# quote-prefix (text starting with "." is treated as ambiguous), which is
# how the source workbook ended up with quotePrefix="1" styled cells.
# Written in this order so the shared-string table is built up the same
# way the original authoring session produced it.
$ws.Range("C3").Value = "'./data/Tony Stark.jpg"
$ws.Range("C4").Value = "'./data/Tom Hank.jpg"
$ws.Range("C2").Value = "'./data/Peter Parker.jpg"

# Widen column C to fit the new path text
$ws.Columns.Item(3).ColumnWidth = 22.6

# Match the author's final selection/cursor position
$ws.Range("J14").Select() | Out-Null
